$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  ,@('Website', 'Projektname', 'Link')
  ,@('Haslehner', 'Elementum Marchtrenk', 'https://www.haslehner.net/projekte/elementum-marchtrenk/')
  ,@('Haslehner', 'Silberweiß 25, Amstetten', 'https://www.haslehner.net/projekte/silberweiss-25-amstetten/')
  ,@('Haslehner', 'rmg56 - Rotenmühlgasse 56, 1120 Wien', 'https://www.haslehner.net/projekte/rmg56-rotenmuehlgasse-56-1120-wien/')
  ,@('Haslehner', 'Ski & Panorama Living - Saalbach-Hinterglemm', 'https://www.haslehner.net/projekte/ski-panorama-living-saalbach-hinterglemm/')
  ,@('Haslehner', 'Hüttelbergstraße 73, 1140 Wien', 'https://www.haslehner.net/projekte/huettelbergstrasse-73-1140-wien/')
  ,@('Haslehner', 'W4 Wallern a.d.T.', 'https://www.haslehner.net/projekte/w4-wallern-adt/')
  ,@('Haslehner', 'Martinstraße 86, 1180 Wien', 'https://www.haslehner.net/projekte/martinstrasse-86-1180-wien/')
  ,@('Haslehner', 'Wohnen 4712 2.0 - Michaelnbach', 'https://www.haslehner.net/projekte/wohnen-4712-20-michaelnbach/')
  ,@('Haslehner', 'Moserstraße 25C, 5020 Salzburg', 'https://www.haslehner.net/projekte/moserstrasse-25c-5020-salzburg/')
  ,@('Haslehner', 'München Stockdorf', 'https://www.haslehner.net/projekte/muenchen-stockdorf/')
  ,@('Haslehner', 'Conrad von Hötzendorf-Platz 3A, 2500 Baden', 'https://www.haslehner.net/projekte/conrad-von-hoetzendorf-platz-3a-2500-baden/')
  ,@('Haslehner', 'Käutzlgasse Salzburg', 'https://www.haslehner.net/projekte/kaeutzlgasse-salzburg/')
  ,@('OOEWohnbau', 'Alkoven, Am Dorfplatz 2 - geförderte Mietwohnungen', 'https://ooewohnbau.at/immobiliensuche/details/alkoven-am-dorfplatz-gefoerderte-wohnungen')
  ,@('OOEWohnbau', 'Andorf, Raaber Straße - 27 geförderte Miet- und Mietkaufwohnungen', 'https://ooewohnbau.at/immobiliensuche/details/wohnungen-andorf-gefoerderte-miet-mietkaufwohnungen')
  ,@('OOEWohnbau', 'Ansfelden, Widistraße - Doppel- und Reihenhäuser / Eigentum ff und Miete mit Kaufoption ff', 'https://ooewohnbau.at/immobiliensuche/details/ansfelden-widistrasse-17-doppel-und-reihenhaeuser-/-eigentum-ff-und-miete-mit-kaufoption-ff')
  ,@('OOEWohnbau', 'Bad Leonfelden, Hochstraße 3 - 15 geförderte Mietwohnungen mit Kaufoption', 'https://ooewohnbau.at/immobiliensuche/details/wohnungen-bad-leonfelden-hochstrasse-gefoerderte-mietkaufwohnungen')
  ,@('OOEWohnbau', 'Grieskirchen, Keplerweg - 21 geförderte Eigentumswohnungen', 'https://ooewohnbau.at/immobiliensuche/details/eigentumswohnungen-in-grieskirchen-21-gefoerderte-eigentumswohnungen-in-grieskirchen')
  ,@('OOEWohnbau', 'Hörsching, Brucknerplatz - Eigentum förderbar oder Miete mit Kaufoption', 'https://ooewohnbau.at/immobiliensuche/details/hoersching-brucknerplatz-28-gefoerderte-eigentumswohnungen')
  ,@('OOEWohnbau', 'Hörsching, Brucknerplatz - Geschäftsfläche Eigentum oder Miete', 'https://ooewohnbau.at/immobiliensuche/details/geschaeftsflaeche-in-hoersching')
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $r = $i + 1
  $rowVals = $data[$i]
  $ws.Cells.Item($r, 1).Value = $rowVals[0]
  $ws.Cells.Item($r, 2).Value = $rowVals[1]
  $ws.Cells.Item($r, 3).Value = $rowVals[2]
}

$headerRange = $ws.Range("A1:C1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

